# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the data refresh captured in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 855
$ws1.Range("F7").Value  = 1231
$ws1.Range("F9").Value  = 808
$ws1.Range("F13").Value = 356
$ws1.Range("F15").Value = 932
$ws1.Range("F16").Value = 9869
$ws1.Range("F17").Value = 616
$ws1.Range("F23").Value = 1759
$ws1.Range("F32").Value = 66
$ws1.Range("F37").Value = 171

# --- Sheet "演出" --------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 128
$ws2.Range("F16").Value = 277

# --- Sheet "本地生活" ----------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 818

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 818
$ws4.Range("F9").Value  = 855
$ws4.Range("F11").Value = 1231
$ws4.Range("F13").Value = 128
$ws4.Range("F14").Value = 808
$ws4.Range("F17").Value = 356
$ws4.Range("F19").Value = 932
$ws4.Range("F20").Value = 9869
$ws4.Range("F22").Value = 616
$ws4.Range("F26").Value = 1759
$ws4.Range("F36").Value = 270
$ws4.Range("F39").Value = 66
$ws4.Range("F47").Value = 171
